$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 38
$ws.Cells.Item(38, 8).Value = 1738.4375  # H38
$ws.Cells.Item(38, 9).Value = 309.77777  # I38
$ws.Cells.Item(38, 10).Value = 3575.2856  # J38
$ws.Cells.Item(38, 11).Value = 929.33331  # K38
$ws.Cells.Item(38, 12).Value = 10725.8568  # L38
$ws.Cells.Item(38, 13).Value = -557.33331  # M38
$ws.Cells.Item(38, 14).Value = -11469.8568  # N38

# ALC row 116
$ws.Cells.Item(116, 8).Value = 27390.25  # H116
$ws.Cells.Item(116, 10).Value = 3185.3333  # J116
$ws.Cells.Item(116, 12).Value = 3185.3333  # L116
$ws.Cells.Item(116, 14).Value = -10069.3333  # N116

# ALC row 137
$ws.Cells.Item(137, 8).Value = 903.86487  # H137
$ws.Cells.Item(137, 9).Value = 888.90625  # I137
$ws.Cells.Item(137, 10).Value = 999.6  # J137
$ws.Cells.Item(137, 11).Value = 2666.71875  # K137
$ws.Cells.Item(137, 12).Value = 2998.8  # L137
$ws.Cells.Item(137, 13).Value = -116.71875  # M137
$ws.Cells.Item(137, 14).Value = -8098.8  # N137

$ws = $wb.Worksheets.Item("ARM")
# ARM row 61
$ws.Cells.Item(61, 8).Value = 189182.19  # H61
$ws.Cells.Item(61, 9).Value = 4240.2974  # I61
$ws.Cells.Item(61, 10).Value = 591702.75  # J61
$ws.Cells.Item(61, 11).Value = 4240.2974  # K61
$ws.Cells.Item(61, 12).Value = 591702.75  # L61
$ws.Cells.Item(61, 13).Value = -4028.2974  # M61
$ws.Cells.Item(61, 14).Value = -592126.75  # N61

# ARM row 132
$ws.Cells.Item(132, 8).Value = 2509.4255  # H132
$ws.Cells.Item(132, 9).Value = 1489.8857  # I132
$ws.Cells.Item(132, 10).Value = 5483.0835  # J132
$ws.Cells.Item(132, 11).Value = 4469.6571  # K132
$ws.Cells.Item(132, 12).Value = 16449.2505  # L132
$ws.Cells.Item(132, 13).Value = -1939.6571  # M132
$ws.Cells.Item(132, 14).Value = -21509.2505  # N132

# ARM row 136
$ws.Cells.Item(136, 8).Value = 189182.19  # H136
$ws.Cells.Item(136, 9).Value = 4240.2974  # I136
$ws.Cells.Item(136, 10).Value = 591702.75  # J136
$ws.Cells.Item(136, 11).Value = 12720.8922  # K136
$ws.Cells.Item(136, 12).Value = 1775108.25  # L136
$ws.Cells.Item(136, 13).Value = -10170.8922  # M136
$ws.Cells.Item(136, 14).Value = -1780208.25  # N136

$ws = $wb.Worksheets.Item("BSM")
# BSM row 134
$ws.Cells.Item(134, 8).Value = 4879.2646  # H134
$ws.Cells.Item(134, 9).Value = 5808.478  # I134
$ws.Cells.Item(134, 10).Value = 2936.3635  # J134
$ws.Cells.Item(134, 11).Value = 17425.434  # K134
$ws.Cells.Item(134, 12).Value = 8809.0905  # L134
$ws.Cells.Item(134, 13).Value = -14890.434  # M134
$ws.Cells.Item(134, 14).Value = -13879.0905  # N134

$ws = $wb.Worksheets.Item("CRP")
# CRP row 15
$ws.Cells.Item(15, 8).Value = 2000  # H15
$ws.Cells.Item(15, 10).Value = 2250  # J15
$ws.Cells.Item(15, 12).Value = 2250  # L15
$ws.Cells.Item(15, 14).Value = -2590  # N15

# CRP row 16
$ws.Cells.Item(16, 8).Value = 2016.2858  # H16
$ws.Cells.Item(16, 9).Value = 1708.9231  # I16
$ws.Cells.Item(16, 10).Value = 2515.75  # J16
$ws.Cells.Item(16, 11).Value = 1708.9231  # K16
$ws.Cells.Item(16, 12).Value = 2515.75  # L16
$ws.Cells.Item(16, 13).Value = -1421.9231  # M16
$ws.Cells.Item(16, 14).Value = -3089.75  # N16

# CRP row 19
$ws.Cells.Item(19, 8).Value = 1933.3334  # H19
$ws.Cells.Item(19, 9).Value = 500  # I19
$ws.Cells.Item(19, 10).Value = 4800  # J19
$ws.Cells.Item(19, 11).Value = 500  # K19
$ws.Cells.Item(19, 12).Value = 4800  # L19
$ws.Cells.Item(19, 13).Value = -330  # M19
$ws.Cells.Item(19, 14).Value = -5140  # N19

# CRP row 24
$ws.Cells.Item(24, 8).Value = 1933.3334  # H24
$ws.Cells.Item(24, 9).Value = 500  # I24
$ws.Cells.Item(24, 10).Value = 4800  # J24
$ws.Cells.Item(24, 11).Value = 500  # K24
$ws.Cells.Item(24, 12).Value = 4800  # L24
$ws.Cells.Item(24, 13).Value = -330  # M24
$ws.Cells.Item(24, 14).Value = -5140  # N24

# CRP row 37
$ws.Cells.Item(37, 8).Value = 257  # H37
$ws.Cells.Item(37, 9).Value = 0  # I37
$ws.Cells.Item(37, 10).Value = 257  # J37
$ws.Cells.Item(37, 11).Value = 0  # K37
$ws.Cells.Item(37, 12).Value = 257  # L37
$ws.Cells.Item(37, 13).ClearContents()  # M37
$ws.Cells.Item(37, 14).Value = -471  # N37

# CRP row 58
$ws.Cells.Item(58, 8).Value = 1491  # H58
$ws.Cells.Item(58, 9).Value = 1289.25  # I58
$ws.Cells.Item(58, 10).Value = 1721.5714  # J58
$ws.Cells.Item(58, 11).Value = 1289.25  # K58
$ws.Cells.Item(58, 12).Value = 1721.5714  # L58
$ws.Cells.Item(58, 13).Value = -1086.25  # M58
$ws.Cells.Item(58, 14).Value = -2127.5714  # N58

# CRP row 113
$ws.Cells.Item(113, 8).Value = 2016.2858  # H113
$ws.Cells.Item(113, 9).Value = 1708.9231  # I113
$ws.Cells.Item(113, 10).Value = 2515.75  # J113
$ws.Cells.Item(113, 11).Value = 1708.9231  # K113
$ws.Cells.Item(113, 12).Value = 2515.75  # L113
$ws.Cells.Item(113, 13).Value = 461.0769  # M113
$ws.Cells.Item(113, 14).Value = -6855.75  # N113

# CRP row 134
$ws.Cells.Item(134, 8).Value = 2645.9348  # H134
$ws.Cells.Item(134, 9).Value = 2543.9023  # I134
$ws.Cells.Item(134, 10).Value = 3482.6  # J134
$ws.Cells.Item(134, 11).Value = 7631.706900000001  # K134
$ws.Cells.Item(134, 12).Value = 10447.8  # L134
$ws.Cells.Item(134, 13).Value = -5096.706900000001  # M134
$ws.Cells.Item(134, 14).Value = -15517.8  # N134

# CRP row 136
$ws.Cells.Item(136, 8).Value = 1491  # H136
$ws.Cells.Item(136, 9).Value = 1289.25  # I136
$ws.Cells.Item(136, 10).Value = 1721.5714  # J136
$ws.Cells.Item(136, 11).Value = 3867.75  # K136
$ws.Cells.Item(136, 12).Value = 5164.7142  # L136
$ws.Cells.Item(136, 13).Value = -1317.75  # M136
$ws.Cells.Item(136, 14).Value = -10264.7142  # N136

$ws = $wb.Worksheets.Item("CUL")
# CUL row 107
$ws.Cells.Item(107, 8).Value = 535.13635  # H107
$ws.Cells.Item(107, 10).Value = 699.2857  # J107
$ws.Cells.Item(107, 12).Value = 2097.8571  # L107
$ws.Cells.Item(107, 14).Value = -5937.8571  # N107

$ws = $wb.Worksheets.Item("GSM")
# GSM row 3
$ws.Cells.Item(3, 8).Value = 4250  # H3
$ws.Cells.Item(3, 10).Value = 7500  # J3
$ws.Cells.Item(3, 12).Value = 7500  # L3
$ws.Cells.Item(3, 14).Value = -7732  # N3

# GSM row 107
$ws.Cells.Item(107, 8).Value = 1191.1852  # H107
$ws.Cells.Item(107, 9).Value = 715.4286  # I107
$ws.Cells.Item(107, 10).Value = 1703.5385  # J107
$ws.Cells.Item(107, 11).Value = 715.4286  # K107
$ws.Cells.Item(107, 12).Value = 1703.5385  # L107
$ws.Cells.Item(107, 13).Value = 1204.5714  # M107
$ws.Cells.Item(107, 14).Value = -5543.538500000001  # N107

# GSM row 132
$ws.Cells.Item(132, 8).Value = 2340.8076  # H132
$ws.Cells.Item(132, 9).Value = 2428.6897  # I132
$ws.Cells.Item(132, 10).Value = 2230  # J132
$ws.Cells.Item(132, 11).Value = 7286.0691  # K132
$ws.Cells.Item(132, 12).Value = 6690  # L132
$ws.Cells.Item(132, 13).Value = -4756.0691  # M132
$ws.Cells.Item(132, 14).Value = -11750  # N132

$ws = $wb.Worksheets.Item("LTW")
# LTW row 29
$ws.Cells.Item(29, 8).Value = 9444.75  # H29
$ws.Cells.Item(29, 10).Value = 9444.75  # J29
$ws.Cells.Item(29, 12).Value = 9444.75  # L29
$ws.Cells.Item(29, 14).Value = -10034.75  # N29

# LTW row 33
$ws.Cells.Item(33, 8).Value = 9415.833000000001  # H33
$ws.Cells.Item(33, 10).Value = 9415.833000000001  # J33
$ws.Cells.Item(33, 12).Value = 9415.833000000001  # L33
$ws.Cells.Item(33, 14).Value = -9995.833000000001  # N33

$ws = $wb.Worksheets.Item("WVR")
# WVR row 3
$ws.Cells.Item(3, 8).Value = 2857034.8  # H3
$ws.Cells.Item(3, 9).Value = 5666733.5  # I3
$ws.Cells.Item(3, 10).Value = 47336  # J3
$ws.Cells.Item(3, 11).Value = 5666733.5  # K3
$ws.Cells.Item(3, 12).Value = 47336  # L3
$ws.Cells.Item(3, 13).Value = -5666619.5  # M3
$ws.Cells.Item(3, 14).Value = -47564  # N3

# WVR row 10
$ws.Cells.Item(10, 8).Value = 3000  # H10
$ws.Cells.Item(10, 10).Value = 3000  # J10
$ws.Cells.Item(10, 12).Value = 3000  # L10
$ws.Cells.Item(10, 14).Value = -3338  # N10

# WVR row 11
$ws.Cells.Item(11, 8).Value = 673350  # H11
$ws.Cells.Item(11, 9).Value = 50  # I11
$ws.Cells.Item(11, 10).Value = 1010000  # J11
$ws.Cells.Item(11, 11).Value = 50  # K11
$ws.Cells.Item(11, 12).Value = 1010000  # L11
$ws.Cells.Item(11, 13).Value = 92  # M11
$ws.Cells.Item(11, 14).Value = -1010284  # N11

# WVR row 22
$ws.Cells.Item(22, 8).Value = 8215  # H22
$ws.Cells.Item(22, 10).Value = 8215  # J22
$ws.Cells.Item(22, 12).Value = 8215  # L22
$ws.Cells.Item(22, 14).Value = -8801  # N22

# WVR row 132
$ws.Cells.Item(132, 8).Value = 1585.0731  # H132
$ws.Cells.Item(132, 9).Value = 1256.2  # I132
$ws.Cells.Item(132, 10).Value = 2482  # J132
$ws.Cells.Item(132, 11).Value = 3768.6  # K132
$ws.Cells.Item(132, 12).Value = 7446  # L132
$ws.Cells.Item(132, 13).Value = -1238.6  # M132
$ws.Cells.Item(132, 14).Value = -12506  # N132
